# Applies the MeetingMinutes.docx edit:
#  - Collapse the spell-check-split runs (w:proofErr wrappers) in the
#    "Members: ..." and "Discussed ... EntityInfo.docx" paragraphs for the
#    2/27 minutes, and the "Members: ..." paragraph for the 3/6 minutes,
#    into single plain runs (text itself is unchanged).
#  - Remove the empty "_GoBack" bookmark paragraph after the 2/27 minutes
#    and add one extra blank paragraph, so the blank-paragraph run between
#    the two meetings grows from 2 (+bookmark paragraph) to 3 plain blanks.
#  - Extend the 3/6 "Discussed" paragraph with the new Spring Break task
#    assignment text, and move the "_GoBack" bookmark to the end of that
#    paragraph.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyInner + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replace paragraph $index (1-based) wholesale -- including its end-of-
# paragraph mark -- with a freshly built paragraph containing $innerXml as
# its content (e.g. one or more <w:r>... runs, bookmarks, etc). This drops
# any stray <w:proofErr/> siblings that Find/Range.Text edits would
# otherwise strand, because the whole paragraph (mark included) is swapped
# out atomically.
function Replace-ParagraphXml($paragraph, [string]$innerXml) {
    $rng = $paragraph.Range
    $xml = New-PkgXml("<w:body><w:p>" + $innerXml + "</w:p></w:body>")
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) 2/27 minutes: "Members: ..." -- merge the spell-check-split runs.
# ---------------------------------------------------------------------
$membersRun1 = '<w:r><w:t>Members: Aaron Raoofi, Savorn Sam, Johnathon Thomas, Emmanuel Ogunkoya, Christopher Kania</w:t></w:r>'
Replace-ParagraphXml $d.Paragraphs(2) $membersRun1

# ---------------------------------------------------------------------
# 2) 2/27 minutes: "Discussed ... EntityInfo.docx" -- merge runs.
# ---------------------------------------------------------------------
$discussed227 = '<w:r><w:t>Discussed which functions/use cases of LMS project fall under the Entity, Boundary and Control classes. Determined a number of mandatory functions in each class type and even added additional functions for Entity and Boundary classes.  The document was approved and finalized by all members present and added to the github (uhdsoftwareengineering/Project-Documentation-Research/) as EntityInfo.docx</w:t></w:r>'
Replace-ParagraphXml $d.Paragraphs(4) $discussed227

# ---------------------------------------------------------------------
# 3) Drop the "_GoBack" bookmark paragraph and insert one extra blank
#    paragraph, so that the 2 blank paragraphs + bookmark paragraph
#    become 3 plain blank paragraphs. Inserting new (bookmark-free)
#    paragraph content right at the start of the bookmark paragraph's
#    own range subsumes/clears that paragraph's bookmark in the same
#    stroke, so no separate delete of a leftover paragraph is needed --
#    paragraph count goes from {bookmark, blank, blank} (3) to
#    {blank, blank, blank} (3), net unchanged.
# ---------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs(5)
$insPoint = $bookmarkPara.Range
$insPoint.Collapse(1)
$insPoint.InsertXML((New-PkgXml("<w:body><w:p/><w:p/></w:body>")))

# ---------------------------------------------------------------------
# 4) 3/6 minutes: "Members: ..." -- merge runs (same pattern as #1).
# ---------------------------------------------------------------------
$membersRun2 = '<w:r><w:t>Members: Aaron Raoofi, Savorn Sam, Johnathon Thomas, Emmanuel Ogunkoya, Christopher Kania</w:t></w:r>'
Replace-ParagraphXml $d.Paragraphs(9) $membersRun2

# ---------------------------------------------------------------------
# 5) 3/6 minutes: "Discussed" paragraph -- append the Spring Break task
#    assignment narrative and move the "_GoBack" bookmark to the end.
# ---------------------------------------------------------------------
$discussedInner = '<w:r><w:t>Discussed</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> the creation of stubs on github for documents that have yet to be </w:t></w:r>' +
    '<w:r><w:t>committed</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
    '<w:r><w:t>Be</w:t></w:r>' +
    '<w:r><w:t>gan to distribute jobs and tasks. Savorn Sam is working on database and use cases. Johnathon Thomas is working on documentation and use cases. Emmanuel Ogunkoya is working on UML diagram. Christopher Kania is working on login page GUI. Aaron Raoofi is working on student view GUI.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Replace-ParagraphXml $d.Paragraphs(11) $discussedInner
